# Restore deployed sources and recovered app.html; prepare local rebuild
#
# 1) Weekly Timesheet: normalise every day to an 8-hour shift at $50/hr,
#    rotate the Bryan/McGill/Hall client assignment, and recompute the
#    SUBTOTAL row (40 reg hrs / $2000).
# 2) Add a new "Jason Schema" worksheet: a flattened/normalised export of
#    the same week (adds Employee, Employee ID and Notes columns).

$wb = $excel.ActiveWorkbook
$ts = $wb.Worksheets.Item(1)

# ---- Weekly Timesheet: update client rotation, hours, and totals ----
# Row 2: 2026-01-12
$ts.Range("B2").Value = "Bryan"
$ts.Range("C2").Value = 8
$ts.Range("F2").Value = 400

# Row 3: 2026-01-13 (client text now resolves to "McGill"; hours/total already correct)
$ts.Range("B3").Value = "McGill"
$ts.Range("C3").Value = 8
$ts.Range("F3").Value = 400

# Row 4: 2026-01-14
$ts.Range("B4").Value = "Hall"
$ts.Range("C4").Value = 8
$ts.Range("F4").Value = 400

# Row 5: 2026-01-15
$ts.Range("B5").Value = "Bryan"
$ts.Range("C5").Value = 8
$ts.Range("F5").Value = 400

# Row 6: 2026-01-16
$ts.Range("B6").Value = "McGill"
$ts.Range("C6").Value = 8
$ts.Range("F6").Value = 400

# SUBTOTAL row (row 8): 40 regular hours, $2000 total
$ts.Range("C8").Value = 40
$ts.Range("D8").Value = "Reg: 40 / OT: 0"
$ts.Range("F8").Value = 2000

# ---- New sheet: "Jason Schema" (flattened export) ----
$js = $wb.Worksheets.Add($null, $ts)
$js.Name = "Jason Schema"

$empName = "Jafid Osorio"
$empId = "emp_AEGgtNcheS2AfGOG"

# Header row (bold)
$js.Range("A1").Value = "Employee"
$js.Range("B1").Value = "Employee ID"
$js.Range("C1").Value = "Date"
$js.Range("D1").Value = "Client"
$js.Range("E1").Value = "Hours"
$js.Range("F1").Value = "Rate"
$js.Range("G1").Value = "Total"
$js.Range("H1").Value = "Type"
$js.Range("I1").Value = "Notes"
$js.Range("A1:I1").Font.Bold = $true
$js.Range("F1:G1").NumberFormat = '"$"#,##0.00'

# Data rows mirror Weekly Timesheet rows 2-6
$dates = @("2026-01-12", "2026-01-13", "2026-01-14", "2026-01-15", "2026-01-16")
$clients = @("Bryan", "McGill", "Hall", "Bryan", "McGill")

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2

    $js.Cells.Item($r, 1).Value = $empName
    $js.Cells.Item($r, 2).Value = $empId

    # Dates are stored as literal text (matching the source export), so
    # force Text format before assignment to stop auto date-parsing.
    $dateCell = $js.Cells.Item($r, 3)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$i]

    $js.Cells.Item($r, 4).Value = $clients[$i]
    $js.Cells.Item($r, 5).Value = 8

    $rateCell = $js.Cells.Item($r, 6)
    $rateCell.NumberFormat = '"$"#,##0.00'
    $rateCell.Value = 50

    $totalCell = $js.Cells.Item($r, 7)
    $totalCell.NumberFormat = '"$"#,##0.00'
    $totalCell.Value = 400

    $js.Cells.Item($r, 8).Value = "Regular"

    # Notes column is blank for every row but the cell still exists.
    $notesCell = $js.Cells.Item($r, 9)
    $notesCell.NumberFormat = "@"
    $notesCell.Value = ""
}

# Column widths to roughly match the source layout
$js.Columns.Item(1).ColumnWidth = 20
$js.Columns.Item(2).ColumnWidth = 18
$js.Columns.Item(3).ColumnWidth = 12
$js.Columns.Item(4).ColumnWidth = 25
$js.Columns.Item(5).ColumnWidth = 8
$js.Columns.Item(6).ColumnWidth = 10
$js.Columns.Item(7).ColumnWidth = 12
$js.Columns.Item(8).ColumnWidth = 10
$js.Columns.Item(9).ColumnWidth = 30

$ts.Activate()
